$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '245.94'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '26.03'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.101'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05597'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.022'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8112'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8458'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1341'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02850'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09375'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001519'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0006040'
$ws.Range('E14').Value = '13OneONEWorstin24h'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.006173'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.562'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06953'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.03210'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1320'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.744'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04693'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001251'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004615'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.00009597'
$ws.Range('E27').Value = '26NitroExNTXBestin24h'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0001390'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03657'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006138'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1054'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002499'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007759'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005327'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1330'
